# Applies the stat corrections to the single-column results table.
# Each row holds one value in its sole cell; we address rows by their
# (1-based) index so duplicate values (e.g. "611") aren't ambiguous.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($rowIndex, $newText) {
    $cell = $t.Rows.Item($rowIndex).Cells.Item(1)
    $cell.Range.Text = $newText
}

Set-CellText 1  "0M"
Set-CellText 2  "0M"
Set-CellText 3  "0M"
Set-CellText 4  "2444"

Set-CellText 6  "0.00080"
Set-CellText 7  "0.00019"
Set-CellText 8  "0.00006"
Set-CellText 9  "0.00026"
Set-CellText 10 "0.00028"
Set-CellText 11 "0.00037"
Set-CellText 12 "0.45251"

Set-CellText 44 "99.99"
Set-CellText 45 "0.45"
Set-CellText 46 "3839"
